$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.030.01'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.910.66'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7883'
$ws.Range("E5").Value = '  +6.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.73'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3161'
$ws.Range("E8").Value = '  +3.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.25'
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06914'
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07981'
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7470'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.908.45'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.227'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.27'
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.063.20'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.02'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.934'
$ws.Range("E18").Value = '  -4.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.78'
$ws.Range("E19").Value = '  +3.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007788'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.899'
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '169.76'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.309'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1383'
$ws.Range("E26").Value = '  +10.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.040'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.381'
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.522'
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.338'
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.112'
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05515'
$ws.Range("E33").Value = '  +4.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.260'
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7365'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.727'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01946'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.796'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.178'
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4436'
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.48'
$ws.Range("E41").Value = '  -0.85%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.897'
$ws.Range("E43").Value = '  -3.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8357'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.56'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.560'
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.809'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '984.83'
$ws.Range("E48").Value = '  +7.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.064.01'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.24'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.497'
$ws.Range("E51").Value = '  +2.63%  '
